# Append a new ticker "200016.IB" as a new row (A40) below the existing
# list of bond tickers (A1:A39), matching the style/format of the
# preceding data cell, and update the active selection to the new cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCell = $ws.Range("A39")
$newCell = $ws.Range("A40")

# Copy the formatting (number format/border/alignment/style) of the last
# existing data cell onto the new cell before setting its value.
$lastCell.Copy()
$newCell.PasteSpecial(-4122)  # xlPasteFormats

$newCell.Value = "200016.IB"

# Mirror the workbook's recorded selection state on the new last cell.
$newCell.Select()
